$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("B1:D14")
$key = $ws.Range("C1:C14")
$range.Sort($key, 2, $null, $null, 1, $null, 1, 1)

$ws.Range("B1:D14").AutoFilter()
